$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 19, shifting old rows 19-21 down to 20-22
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with its data
$ws.Cells.Item(19, 1).Value = 4
$ws.Cells.Item(19, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(19, 3).Value = "Los Lagos"
$ws.Cells.Item(19, 4).Value = 44876
$ws.Cells.Item(19, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(19, 5).Value = 10
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100107
$ws.Cells.Item(19, 8).Value = "Otros"
$ws.Cells.Item(19, 9).Value = 100107002
$ws.Cells.Item(19, 10).Value = "Chirimoya"
$ws.Cells.Item(19, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 300
$ws.Cells.Item(19, 14).Value = 22000
$ws.Cells.Item(19, 15).Value = 22500
$ws.Cells.Item(19, 16).Value = 22250
$ws.Cells.Item(19, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(19, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 19).Value = 2781
$ws.Cells.Item(19, 20).Value = 8

Write-Output "Done"
